# Trading update: 2026-02-17 13:28:20
# Append the newly-opened MarketMaking trade (#48) as row 49 on both the
# "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $row = 49

    $ws.Cells.Item($row, 1).Value = 48

    # Dates like "2026-02-17" are auto-parsed into Excel date serials, so
    # force the cell to text first and then clear the resulting formatting
    # so the stored value is the literal string (matching the rest of the
    # column, which holds plain text dates).
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "13:28:19"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.95

    # Trade is still OPEN, so there is no exit price yet.
    $ws.Cells.Item($row, 7).Value = " "

    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 97.51239312960779
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # No exit reason yet either, since the trade has not closed.
    $ws.Cells.Item($row, 16).Value = " "

    $ws.Cells.Item($row, 17).Value = 0
}
